# Updates cryptos list data (prices & 1h volume change) per the commit diff.
# D-column values are numeric-looking text (e.g. "34.953.11") that must stay
# text -- a leading apostrophe forces Excel to keep them as text while
# preserving the original "General" number format (matches original style).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'34.953.11"
$ws.Range("E2").Value = "  +1.01%  "

$ws.Range("D3").Value = "'1.846.05"
$ws.Range("E3").Value = "  +1.88%  "

$ws.Range("E4").Value = "  +0.08%  "

$ws.Range("D5").Value = "'227.75"
$ws.Range("E5").Value = "  +0.77%  "

$ws.Range("D6").Value = "'0.613"
$ws.Range("E6").Value = "  +1.93%  "

$ws.Range("D7").Value = "'0.999"
$ws.Range("E7").Value = "  +0.03%  "

$ws.Range("D8").Value = "'42.11"
$ws.Range("E8").Value = "  +14.77%  "

$ws.Range("D9").Value = "'0.305"
$ws.Range("E9").Value = "  +3.82%  "

$ws.Range("D10").Value = "'0.0693"
$ws.Range("E10").Value = "  +1.26%  "

$ws.Range("E11").Value = "  +3.42%  "

$ws.Range("D12").Value = "'2.113.81"
$ws.Range("E12").Value = "  +1.94%  "

$ws.Range("D13").Value = "'11.66"
$ws.Range("E13").Value = "  +2.80%  "

$ws.Range("D14").Value = "'1.842.72"
$ws.Range("E14").Value = "  +1.73%  "

$ws.Range("D15").Value = "'4.76"
$ws.Range("E15").Value = "  +7.07%  "

$ws.Range("E16").Value = "  +4.28%  "

$ws.Range("D17").Value = "'34.876.77"
$ws.Range("E17").Value = "  +0.91%  "

$ws.Range("D18").Value = "'69.47"
$ws.Range("E18").Value = "  +1.02%  "

$ws.Range("D19").Value = "'244.81"
$ws.Range("E19").Value = "  +0.58%  "

$ws.Range("D20").Value = "'0.0₃0790"
$ws.Range("E20").Value = "  +1.33%  "

$ws.Range("D21").Value = "'12.12"
$ws.Range("E21").Value = "  +7.88%  "

$ws.Range("D22").Value = "'4.75"
$ws.Range("E22").Value = "  +15.24%  "

$ws.Range("E23").Value = "  +0.17%  "

$ws.Range("E24").Value = "  -1.12%  "

$ws.Range("D25").Value = "'171.97"
$ws.Range("E25").Value = "  +0.07%  "

$ws.Range("D26").Value = "'7.96"
$ws.Range("E26").Value = "  +1.13%  "

$ws.Range("D27").Value = "'17.85"
$ws.Range("E27").Value = "  +3.41%  "

$ws.Range("D28").Value = "'0.122"
$ws.Range("E28").Value = "  +1.06%  "

$ws.Range("E29").Value = "  +0.04%  "

$ws.Range("E30").Value = "  +8.31%  "

$ws.Range("D31").Value = "'3.98"
$ws.Range("E31").Value = "  +3.71%  "

$ws.Range("D32").Value = "'4.03"
$ws.Range("E32").Value = "  +2.65%  "

$ws.Range("E33").Value = "  +3.51%  "

$ws.Range("D34").Value = "'1.92"
$ws.Range("E34").Value = "  +5.64%  "

$ws.Range("D35").Value = "'91.04"
$ws.Range("E35").Value = "  +11.59%  "

$ws.Range("E36").Value = "  +2.51%  "

$ws.Range("B37").Value = "TrustWalletToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D37").Value = "'1.10"
$ws.Range("E37").Value = "  +2.81%  "

$ws.Range("B38").Value = "Maker"
$ws.Range("C38").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D38").Value = "'1.345.25"
$ws.Range("E38").Value = "  -1.53%  "

$ws.Range("D39").Value = "'1.03"
$ws.Range("E39").Value = "  +9.33%  "

$ws.Range("E40").Value = "  +2.76%  "

$ws.Range("E41").Value = "  +3.54%  "

$ws.Range("D42").Value = "'14.88"
$ws.Range("E42").Value = "  +8.59%  "

$ws.Range("E43").Value = "  +6.86%  "

$ws.Range("E44").Value = "  +1.81%  "

$ws.Range("E46").Value = "  +2.96%  "

$ws.Range("B47").Value = "FraxShare"
$ws.Range("C47").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D47").Value = "'6.05"
$ws.Range("E47").Value = "  +3.71%  "

$ws.Range("B48").Value = "RocketPoolETH"
$ws.Range("C48").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D48").Value = "'2.012.27"
$ws.Range("E48").Value = "  +1.96%  "

$ws.Range("E49").Value = "  +0.11%  "

$ws.Range("D50").Value = "'103.34"
$ws.Range("E50").Value = "  +0.33%  "

$ws.Range("D51").Value = "'7.25"
$ws.Range("E51").Value = "  +5.09%  "
